$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold + border) from F1 into the new header cells G1:L1
$ws.Range("F1").Copy()
$ws.Range("G1:L1").PasteSpecial(-4122)

# New header labels
$ws.Range("G1").Value = "W1"
$ws.Range("H1").Value = "W1 Reason"
$ws.Range("I1").Value = "W2"
$ws.Range("J1").Value = "W2 Reason"
$ws.Range("K1").Value = "W3"
$ws.Range("L1").Value = "W3 Reason"

# Weekly absence data (week number + reason) for each student row
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = "Dentist"
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = "Dentist"
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = "Dentist"
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = "Ill"
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = "Ill"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = "Doctors"
$ws.Range("K4").Value = 7
$ws.Range("L4").Value = "Doctors"
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = "Trip"
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = "Trip"
$ws.Range("K6").Value = 7
$ws.Range("L6").Value = "Trip"
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = "Not Well"
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = "Holiday"
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = "Holiday"
$ws.Range("K9").Value = 7
$ws.Range("L9").Value = "Ill"
$ws.Range("K10").Value = 7
$ws.Range("L10").Value = "Holiday"
$ws.Range("K11").Value = 7
$ws.Range("L11").Value = "Unauthorised"
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = "Doctors"
$ws.Range("K13").Value = 7
$ws.Range("L13").Value = "Not Well"
$ws.Range("K14").Value = 7
$ws.Range("L14").Value = "Ill"
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = "Doctors"
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = "Not Well"
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = "Ill"
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = "Ill"
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = "Doctors"
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = "Doctors"
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = "Unauthorised"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = "Unauthorised"
